$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------------------
# 1) "All: Good communication and Clean code" -> "All: Good communication and
#    clean code", re-split across three runs (no more proofErr wrapper).
# ---------------------------------------------------------------------------
$allPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith('All:')) {
        $allPara = $p
        break
    }
}
if ($allPara -eq $null) {
    throw "Could not find the 'All:' paragraph"
}
$allXml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:r><w:t>All:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Good communication and cl</w:t></w:r>' +
    '<w:r><w:t>ean code</w:t></w:r>' +
    '</w:p>'
$allPara.Range.InsertXML($allXml)

# ---------------------------------------------------------------------------
# 2) "Schedule:" paragraph -> "Schedule: " (single run, trailing space,
#    bookmark removed).
# ---------------------------------------------------------------------------
$schedulePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith('Schedule:')) {
        $schedulePara = $p
        break
    }
}
if ($schedulePara -eq $null) {
    throw "Could not find the 'Schedule:' paragraph"
}
$scheduleXml = '<w:p xmlns:w="' + $wNs + '"><w:r><w:t xml:space="preserve">Schedule: </w:t></w:r></w:p>'
$schedulePara.Range.InsertXML($scheduleXml)

# ---------------------------------------------------------------------------
# 3) Insert the five new schedule-detail paragraphs right after "Schedule: ".
# ---------------------------------------------------------------------------
$endOfSchedule = $d.Range($schedulePara.Range.End, $schedulePara.Range.End)

$nd = [string][char]0x2026

$newParasXml =
    '<w:p xmlns:w="' + $wNs + '"><w:r><w:t>3-31: Group meeting</w:t></w:r></w:p>' +
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:t>4</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:t>-8: Half of our coding done</w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:t>4-9:2</w:t></w:r>' +
        '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>nd</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> Group </w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>Meeting</w:t></w:r>' +
        '<w:r><w:t>(</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t>If Possible)</w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:t>4-15</w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>:Finish</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t xml:space="preserve"> basic coding</w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:t>4-16-</w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>' + $nd + ' :Meetings</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t>/Cleanup</w:t></w:r>' +
    '</w:p>'

$endOfSchedule.InsertXML($newParasXml)

Write-Output "done"
